$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text (matches source inlineStr cells)
$textCells = @("D5", "D6", "D9", "D11", "D12", "D20", "D21", "D24", "D25", "D28", "D29", "D30", "D31", "D32", "D37", "D39", "D40", "D41", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates row by row
# Row 2
$ws.Range("D2").Value = "63.193.99"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3
$ws.Range("D3").Value = "3.153.72"
$ws.Range("E3").Value = "  +1.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "589.70"
$ws.Range("E5").Value = "  -1.92%  "

# Row 6
$ws.Range("D6").Value = "138.37"
$ws.Range("E6").Value = "  -3.09%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "3.153.11"
$ws.Range("E8").Value = "  +1.24%  "

# Row 9
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  -0.24%  "

# Row 10
$ws.Range("E10").Value = "  -1.30%  "

# Row 11
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  -0.46%  "

# Row 12
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -1.30%  "

# Row 13
$ws.Range("E13").Value = "  -2.15%  "

# Row 14
$ws.Range("E14").Value = "  -2.24%  "

# Row 15
$ws.Range("D15").Value = "3.675.80"
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("E16").Value = "  +1.21%  "

# Row 17
$ws.Range("D17").Value = "3.151.25"
$ws.Range("E17").Value = "  +0.73%  "

# Row 18
$ws.Range("D18").Value = "63.191.40"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").Value = "477.36"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  -3.55%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("E23").Value = "  +1.99%  "

# Row 24
$ws.Range("D24").Value = "84.74"
$ws.Range("E24").Value = "  -3.16%  "

# Row 25
$ws.Range("D25").Value = "13.02"
$ws.Range("E25").Value = "  -1.98%  "

# Row 26
$ws.Range("E26").Value = "  +0.03%  "

# Row 27
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("D28").Value = "7.13"
$ws.Range("E28").Value = "  +1.49%  "

# Row 29
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -3.04%  "

# Row 30
$ws.Range("D30").Value = "2.12"
$ws.Range("E30").Value = "  +3.68%  "

# Row 31
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.08%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "26.92"
$ws.Range("E32").Value = "  -0.74%  "

# Row 33
$ws.Range("E33").Value = "  -4.10%  "

# Row 34
$ws.Range("E34").Value = "  -3.30%  "

# Row 35
$ws.Range("E35").Value = "  -2.34%  "

# Row 36
$ws.Range("E36").Value = "  -2.71%  "

# Row 37
$ws.Range("D37").Value = "52.50"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -6.47%  "

# Row 39
$ws.Range("D39").Value = "0.0389"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40
$ws.Range("D40").Value = "416.44"
$ws.Range("E40").Value = "  -4.42%  "

# Row 41
$ws.Range("D41").Value = "2.77"
$ws.Range("E41").Value = "  -5.70%  "

# Row 42
$ws.Range("E42").Value = "  +0.79%  "

# Row 43
$ws.Range("D43").Value = "2.930.31"
$ws.Range("E43").Value = "  +2.64%  "

# Row 44
$ws.Range("E44").Value = "  -6.09%  "

# Row 45
$ws.Range("D45").Value = "0.263"
$ws.Range("E45").Value = "  +1.87%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.13"
$ws.Range("E47").Value = "  -2.91%  "

# Row 48
$ws.Range("D48").Value = "25.47"
$ws.Range("E48").Value = "  -0.89%  "

# Row 49
$ws.Range("E49").Value = "  +0.35%  "

# Row 50
$ws.Range("E50").Value = "  -7.24%  "

# Row 51
$ws.Range("D51").Value = "121.34"
$ws.Range("E51").Value = "  -0.15%  "
